$d = $word.ActiveDocument

# 1) Replace the duplicated "double switch" sentence in the first paragraph
#    with the three new sentences about the amperemeter/conductor/resistor.
$d.Content.Find.Execute(
    "На схему был добавлен двойной переключатель. На схему был добавлен двойной переключатель. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "На схему был добавлен амперметр. На схему был добавлен проводник. На схему был добавлен резистор. ",
    2
)

# 2) Remove the trailing paragraph that used to read
#    "На схему был добавлен амперметр. " (its whole paragraph, mark included,
#    is dropped now that the sentence was folded into the first paragraph).
$paraCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($paraCount)
$lastPara.Range.Delete()
